# QuestionBank-Template-V2.xlsx — add a new multiple-choice question row
# and switch the active sheet/selection back to "MultipleChoice".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MultipleChoice")

# --- Append the new question as row 18 (A:I = question,a,b,c,d,e,ans,image,category)
$ws.Range("A18").Value = "How many corners are there in square?"
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = "d"
$ws.Range("H18").Value = "square.png"
$ws.Range("I18").Value = "maths"

# --- New "image" column (H) gets a best-fit-style width for the new text
$ws.Columns.Item(8).ColumnWidth = 9.67

# --- Restore selection on the new row and make MultipleChoice the active/visible tab
$ws.Range("I18").Select() | Out-Null
$ws.Activate() | Out-Null
